$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1244.2222
$ws.Range("I2").Value = 440
$ws.Range("K2").Value = 440
$ws.Range("M2").Value = -327
$ws.Range("H9").Value = 2206.1904
$ws.Range("I9").Value = 2623.0588
$ws.Range("J9").Value = 434.5
$ws.Range("K9").Value = 2623.0588
$ws.Range("L9").Value = 434.5
$ws.Range("M9").Value = -2454.0588
$ws.Range("N9").Value = -772.5
$ws.Range("H17").Value = 1734.7715
$ws.Range("J17").Value = 1734.7715
$ws.Range("L17").Value = 5204.3145
$ws.Range("N17").Value = -5540.3145
$ws.Range("H18").Value = 2120
$ws.Range("I18").Value = 1640
$ws.Range("K18").Value = 1640
$ws.Range("M18").Value = -1356
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H33").Value = 256.4
$ws.Range("I33").Value = 292.16666
$ws.Range("K33").Value = 292.16666
$ws.Range("M33").Value = -63.16665999999998
$ws.Range("H40").Value = 62516256
$ws.Range("I40").Value = 2996.25
$ws.Range("K40").Value = 2996.25
$ws.Range("M40").Value = -2821.25
$ws.Range("H112").Value = 3998.6843
$ws.Range("J112").Value = 4165.278
$ws.Range("L112").Value = 12495.834
$ws.Range("N112").Value = -14711.834
$ws.Range("H129").Value = 4050
$ws.Range("I129").Value = 823.1429000000001
$ws.Range("J129").Value = 11579.333
$ws.Range("K129").Value = 2469.4287
$ws.Range("L129").Value = 34737.999
$ws.Range("M129").Value = 2530.5713
$ws.Range("N129").Value = -44737.999
$ws.Range("H132").Value = 2599.5264
$ws.Range("I132").Value = 2381.8823
$ws.Range("J132").Value = 4449.5
$ws.Range("K132").Value = 7145.646900000001
$ws.Range("L132").Value = 13348.5
$ws.Range("M132").Value = -4615.646900000001
$ws.Range("N132").Value = -18408.5
$ws.Range("H138").Value = 3418.6072
$ws.Range("I138").Value = 2234.8262
$ws.Range("J138").Value = 4243.6665
$ws.Range("K138").Value = 6704.4786
$ws.Range("L138").Value = 12730.9995
$ws.Range("M138").Value = -1564.4786
$ws.Range("N138").Value = -23010.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9699.950000000001
$ws.Range("I32").Value = 8839.945
$ws.Range("K32").Value = 8839.945
$ws.Range("M32").Value = -8552.945
$ws.Range("H61").Value = 3966847
$ws.Range("I61").Value = 4882550.5
$ws.Range("J61").Value = 838193.8
$ws.Range("K61").Value = 4882550.5
$ws.Range("L61").Value = 838193.8
$ws.Range("M61").Value = -4882338.5
$ws.Range("N61").Value = -838617.8
$ws.Range("H92").Value = 52550
$ws.Range("J92").Value = 52550
$ws.Range("L92").Value = 52550
$ws.Range("N92").Value = -57542
$ws.Range("H96").Value = 44895.75
$ws.Range("J96").Value = 44895.75
$ws.Range("L96").Value = 44895.75
$ws.Range("N96").Value = -50387.75
$ws.Range("I132").Value = 1990.0698
$ws.Range("J132").Value = 25002904
$ws.Range("K132").Value = 5970.2094
$ws.Range("L132").Value = 75008712
$ws.Range("M132").Value = -3440.2094
$ws.Range("N132").Value = -75013772
$ws.Range("H133").Value = 94000.5
$ws.Range("I133").Value = 94000.5
$ws.Range("K133").Value = 94000.5
$ws.Range("M133").Value = -91470.5
$ws.Range("H136").Value = 3966847
$ws.Range("I136").Value = 4882550.5
$ws.Range("J136").Value = 838193.8
$ws.Range("K136").Value = 14647651.5
$ws.Range("L136").Value = 2514581.4
$ws.Range("M136").Value = -14645101.5
$ws.Range("N136").Value = -2519681.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4351494
$ws.Range("I134").Value = 3686.6843
$ws.Range("K134").Value = 11060.0529
$ws.Range("M134").Value = -8525.052899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29415250
$ws.Range("I31").Value = 50002964
$ws.Range("K31").Value = 50002964
$ws.Range("M31").Value = -50002669
$ws.Range("H34").Value = 29415250
$ws.Range("I34").Value = 50002964
$ws.Range("K34").Value = 50002964
$ws.Range("M34").Value = -50002762
$ws.Range("H58").Value = 2323.1924
$ws.Range("I58").Value = 1720.2632
$ws.Range("J58").Value = 3959.7144
$ws.Range("K58").Value = 1720.2632
$ws.Range("L58").Value = 3959.7144
$ws.Range("M58").Value = -1517.2632
$ws.Range("N58").Value = -4365.7144
$ws.Range("H134").Value = 2215.5334
$ws.Range("I134").Value = 2011.6818
$ws.Range("J134").Value = 2776.125
$ws.Range("K134").Value = 6035.0454
$ws.Range("L134").Value = 8328.375
$ws.Range("M134").Value = -3500.0454
$ws.Range("N134").Value = -13398.375
$ws.Range("H136").Value = 2323.1924
$ws.Range("I136").Value = 1720.2632
$ws.Range("J136").Value = 3959.7144
$ws.Range("K136").Value = 5160.7896
$ws.Range("L136").Value = 11879.1432
$ws.Range("M136").Value = -2610.7896
$ws.Range("N136").Value = -16979.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 150.16667
$ws.Range("I2").Value = 80
$ws.Range("K2").Value = 480
$ws.Range("M2").Value = -367
$ws.Range("H11").Value = 111116770
$ws.Range("I11").Value = 142858080
$ws.Range("J11").Value = 22166
$ws.Range("K11").Value = 428574240
$ws.Range("L11").Value = 66498
$ws.Range("M11").Value = -428574100
$ws.Range("N11").Value = -66778
$ws.Range("H26").Value = 6452
$ws.Range("I26").Value = 1075.8
$ws.Range("K26").Value = 3227.4
$ws.Range("M26").Value = -2939.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4404.5356
$ws.Range("I126").Value = 4121.65
$ws.Range("K126").Value = 12364.95
$ws.Range("M126").Value = -9894.949999999999
$ws.Range("H132").Value = 1757724.9
$ws.Range("I132").Value = 3305.796
$ws.Range("K132").Value = 9917.387999999999
$ws.Range("M132").Value = -7387.387999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 4987.6665
$ws.Range("I23").Value = 4987.6665
$ws.Range("K23").Value = 4987.6665
$ws.Range("M23").Value = -4757.6665
$ws.Range("H40").Value = 5745.25
$ws.Range("I40").Value = 4719.6924
$ws.Range("K40").Value = 4719.6924
$ws.Range("M40").Value = -4583.6924
$ws.Range("H136").Value = 3899.5518
$ws.Range("I136").Value = 2320.6316
$ws.Range("K136").Value = 6961.8948
$ws.Range("M136").Value = -4411.8948

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 68130
$ws.Range("J118").Value = 68130
$ws.Range("L118").Value = 68130
$ws.Range("N118").Value = -71444
$ws.Range("H132").Value = 516607.78
$ws.Range("I132").Value = 4031.5518
$ws.Range("J132").Value = 2003078.8
$ws.Range("K132").Value = 12094.6554
$ws.Range("L132").Value = 6009236.4
$ws.Range("M132").Value = -9564.6554
$ws.Range("N132").Value = -6014296.4
$ws.Range("H136").Value = 324615.16
$ws.Range("I136").Value = 2244.524
$ws.Range("K136").Value = 6733.572
$ws.Range("M136").Value = -4183.572
$ws.Range("H140").Value = 52910.57
$ws.Range("J140").Value = 52910.57
$ws.Range("L140").Value = 52910.57
$ws.Range("N140").Value = -63270.57
$ws.Range("H141").Value = 89960.11
$ws.Range("J141").Value = 89960.11
$ws.Range("L141").Value = 89960.11
$ws.Range("N141").Value = -100320.11
